$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$full = "Assignment 1 due on September 26`rAutograder out soon`r`rTCP is provides stream-level guarantees`rNot message-level`r"
$tr.Text = $full

$after = $shape.TextFrame.TextRange
$p2 = $after.Paragraphs(2,1)
$p2.Font.Bold = $false
Write-Output "p2 len = $($p2.Length)"
$run1 = $p2.Characters(1,10)
Write-Output "run1=[$($run1.Text)]"
$run2 = $p2.Characters(11,9)
Write-Output "run2=[$($run2.Text)]"
